$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of ISIN data (row 10), matching the layout/format of the existing
# rows above it (series_mosb, series_short, isin, maturity date).

# Column C (isin) - plain text, General style (same as other C column cells)
$ws.Range("C10").Value = "PIBD0726B627"

# Column A (series_mosb) - copy format from A6 (mmm-yy numfmt + quotePrefix),
# then set the text value (leading apostrophe keeps it text instead of a date)
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "'7-62"

# Column B (series_short) - copy format from B6 (quotePrefix/general),
# then set the text value
$ws.Range("B6").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "'762"

# Column D (maturity) - copy format from D6 (custom date numfmt),
# then set the numeric (date serial) value
$ws.Range("D6").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = 46067

# Update the active selection to match where the user left off editing
$ws.Range("D14").Select()
